{"js": "// Replace the 25 \"two-digit number \u00f7 one-digit number =\" expressions in\n// the practice-sheet table with newly generated problems, in document\n// order. Source values are not all unique (e.g. \"60\u00f74=\" occurs twice,\n// mapped to two different replacements), so each entry tracks which\n// occurrence (0-based) of its source text it targets.\nconst replacements = [\n  { text: \"60\u00f74=\", occurrence: 0, replacement: \"81\u00f75=\" },\n  { text: \"37\u00f75=\", occurrence: 0, replacement: \"31\u00f76=\" },\n  { text: \"38\u00f76=\", occurrence: 0, replacement: \"81\u00f76=\" },\n  { text: \"56\u00f77=\", occurrence: 0, replacement: \"17\u00f76=\" },\n  { text: \"93\u00f79=\", occurrence: 0, replacement: \"12\u00f77=\" },\n  { text: \"35\u00f76=\", occurrence: 0, replacement: \"80\u00f76=\" },\n  { text: \"28\u00f74=\", occurrence: 0, replacement: \"18\u00f76=\" },\n  { text: \"51\u00f74=\", occurrence: 0, replacement: \"39\u00f73=\" },\n  { text: \"12\u00f72=\", occurrence: 0, replacement: \"36\u00f73=\" },\n  { text: \"89\u00f73=\", occurrence: 0, replacement: \"11\u00f75=\" },\n  { text: \"54\u00f74=\", occurrence: 0, replacement: \"20\u00f72=\" },\n  { text: \"92\u00f73=\", occurrence: 0, replacement: \"58\u00f76=\" },\n  { text: \"71\u00f72=\", occurrence: 0, replacement: \"26\u00f72=\" },\n  { text: \"37\u00f72=\", occurrence: 0, replacement: \"16\u00f77=\" },\n  { text: \"26\u00f73=\", occurrence: 0, replacement: \"80\u00f73=\" },\n  { text: \"98\u00f76=\", occurrence: 0, replacement: \"26\u00f78=\" },\n  { text: \"96\u00f78=\", occurrence: 0, replacement: \"71\u00f77=\" },\n  { text: \"45\u00f76=\", occurrence: 0, replacement: \"61\u00f73=\" },\n  { text: \"52\u00f73=\", occurrence: 0, replacement: \"69\u00f79=\" },\n  { text: \"31\u00f74=\", occurrence: 0, replacement: \"65\u00f78=\" },\n  { text: \"60\u00f74=\", occurrence: 1, replacement: \"65\u00f73=\" },\n  { text: \"76\u00f78=\", occurrence: 0, replacement: \"25\u00f75=\" },\n  { text: \"37\u00f79=\", occurrence: 0, replacement: \"64\u00f72=\" },\n  { text: \"35\u00f73=\", occurrence: 0, replacement: \"56\u00f77=\" },\n  { text: \"10\u00f73=\", occurrence: 0, replacement: \"43\u00f74=\" },\n];\n\n// Group by source text so we only issue one search per distinct string.\nconst byText = new Map();\nfor (const item of replacements) {\n  if (!byText.has(item.text)) {\n    byText.set(item.text, []);\n  }\n  byText.get(item.text).push(item);\n}\n\nconst searchResultsByText = new Map();\nfor (const text of byText.keys()) {\n  const results = context.document.body.search(text, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  searchResultsByText.set(text, results);\n}\nawait context.sync();\n\nfor (const [text, items] of byText.entries()) {\n  const results = searchResultsByText.get(text);\n  for (const item of items) {\n    const range = results.items[item.occurrence];\n    range.insertText(item.replacement, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Replace the 25 \"two-digit number \u00f7 one-digit number =\" expressions in\n# the practice-sheet table with newly generated problems, in document\n# order. Source values are not all unique (e.g. \"60\u00f74=\" occurs twice,\n# mapped to two different replacements), so matches are located by\n# sweeping the document once, left to right, and pairing each match (in\n# order) with its corresponding replacement below.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    \"81\u00f75=\",\n    \"31\u00f76=\",\n    \"81\u00f76=\",\n    \"17\u00f76=\",\n    \"12\u00f77=\",\n    \"80\u00f76=\",\n    \"18\u00f76=\",\n    \"39\u00f73=\",\n    \"36\u00f73=\",\n    \"11\u00f75=\",\n    \"20\u00f72=\",\n    \"58\u00f76=\",\n    \"26\u00f72=\",\n    \"16\u00f77=\",\n    \"80\u00f73=\",\n    \"26\u00f78=\",\n    \"71\u00f77=\",\n    \"61\u00f73=\",\n    \"69\u00f79=\",\n    \"65\u00f78=\",\n    \"65\u00f73=\",\n    \"25\u00f75=\",\n    \"64\u00f72=\",\n    \"56\u00f77=\",\n    \"43\u00f74=\"\n)\n\n$searchTexts = @(\n    \"60\u00f74=\", \"37\u00f75=\", \"38\u00f76=\", \"56\u00f77=\", \"93\u00f79=\",\n    \"35\u00f76=\", \"28\u00f74=\", \"51\u00f74=\", \"12\u00f72=\", \"89\u00f73=\",\n    \"54\u00f74=\", \"92\u00f73=\", \"71\u00f72=\", \"37\u00f72=\", \"26\u00f73=\",\n    \"98\u00f76=\", \"96\u00f78=\", \"45\u00f76=\", \"52\u00f73=\", \"31\u00f74=\",\n    \"60\u00f74=\", \"76\u00f78=\", \"37\u00f79=\", \"35\u00f73=\", \"10\u00f73=\"\n)\n\n# Sweep left to right, one match at a time, capturing a Duplicate (a\n# detached snapshot range) of each expression cell's text run in\n# document order. This resolves each of the 25 target cells exactly\n# once against the ORIGINAL text, even though some search strings\n# (e.g. \"60\u00f74=\") repeat.\n$matchRanges = New-Object System.Collections.ArrayList\n$cursor = $d.Content\n$cursor.Start = 0\nfor ($idx = 0; $idx -lt $searchTexts.Count; $idx++) {\n    $probe = $d.Content\n    $probe.Start = $cursor.Start\n    $probe.End = $d.Content.End\n    $found = $probe.Find.Execute($searchTexts[$idx], $false, $false, $false, $false, $false, $true, 1, $false, $null, 0)\n    [void]$matchRanges.Add($probe.Duplicate)\n    $cursor.Start = $probe.End\n}\n\n# Apply the replacements from the last match to the first so that\n# earlier (not-yet-applied) ranges keep their original offsets.\nfor ($idx = $matchRanges.Count - 1; $idx -ge 0; $idx--) {\n    $matchRanges[$idx].Text = $replacements[$idx]\n}\n"}
